# Timesheet update: the "01 Sept 24" time entry for yeasir afgan is removed,
# and the "29 Aug 24" entry's Start/End times (and therefore Hours Worked,
# Total, and Grand Total) are revised to 13:48 - 17:48 / 4.00 hours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the second time-entry row (originally row 3: "01 Sept 24" / 20:32-03:35).
# This shifts the old Total row (4) and Grand Total row (5) up to rows 3 and 4.
$ws.Rows.Item(3).Delete()

# Revise the remaining time entry (row 2): new Start Time, End Time, Hours Worked.
$ws.Range("C2").Value = "13:48"
$ws.Range("D2").Value = "17:48"
$ws.Range("E2").Value = 4

# Refresh the Total / Grand Total cells (now rows 3 & 4) to "4.00". They are
# stored as literal text in the sheet (like the original "14.56"), so write
# them via a text formula and then paste-special-values: a direct
# Range.Value assignment of a numeric-looking string like "4.00" gets
# auto-coerced to the number 4 by Excel, which would also require a new
# Text number-format (and a styles.xml change that the original edit didn't
# make). Formula + paste-values keeps the literal string without touching
# any cell's number format/style.
$ws.Range("E3").Formula = "=""4.00"""
$ws.Range("E4").Formula = "=""4.00"""
$ws.Range("E3:E4").Copy()
$ws.Range("E3:E4").PasteSpecial(-4163)

# Column B narrows from width 12 to width 11 (OOXML units); Excel's
# ColumnWidth property is offset from the raw OOXML column width, so 11
# raw units corresponds to a ColumnWidth of 10.17 here (matches column A/D/E,
# whose raw widths of 14/10/14 read back as ColumnWidth 13.17/9.17/13.17).
$ws.Columns.Item(2).ColumnWidth = 10.17
